# [ADDITIONAL SCRAPING] added code to scrape more data about a player's
# batting performance in a match, also updated the excel sheets.
#
#  1. Insert a new "Player Info" worksheet as the FIRST sheet, containing
#     ID / NAME / BATTING_HAND / BOWL_STYLE for player 5861.
#  2. On the "ODI Batting" sheet, rename the MATCH_CARD_LINK column to
#     MATCH_CODE and replace the full scorecard URL values with just the
#     bare match-code number.
#  3. Do the same (MATCH_CARD_LINK -> MATCH_CODE) on the "ODI Bowling"
#     sheet.

$wb = $excel.ActiveWorkbook

# --- 1. Add the new "Player Info" sheet, positioned first -------------------

$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "Player Info"

$playerInfo = $wb.Worksheets.Item("Player Info")
$battingSheet = $wb.Worksheets.Item("ODI Batting")

# Reuse the bold/bordered header style already used by the other sheets'
# header rows (style index stays shared instead of minting a new one).
$battingSheet.Range("A1").Copy()
$playerInfo.Range("A1:D1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$playerInfo.Range("A1").Value = "ID"
$playerInfo.Range("B1").Value = "NAME"
$playerInfo.Range("C1").Value = "BATTING_HAND"
$playerInfo.Range("D1").Value = "BOWL_STYLE"

# Player id looks numeric but should be stored as text (like the rest of
# the scraped data); force text formatting, write it, then drop back to
# the default (unstyled) cell format used by the rest of the data rows.
$defaultStyle = $playerInfo.Range("B2").Style
$playerInfo.Range("A2").NumberFormat = "@"
$playerInfo.Range("A2").Value = "5861"
$playerInfo.Range("A2").Style = $defaultStyle

$playerInfo.Range("B2").Value = "Shahbaz Ahmed"
$playerInfo.Range("C2").Value = "Left Handed"
$playerInfo.Range("D2").Value = "Left Arm Orthodox"

$playerInfo.Range("A1").Select() | Out-Null

# --- 2. ODI Batting sheet: MATCH_CARD_LINK -> MATCH_CODE --------------------

$battingSheet.Range("D1").Value = "MATCH_CODE"

$battingDefaultStyle = $battingSheet.Range("A2").Style

$battingSheet.Range("D2").NumberFormat = "@"
$battingSheet.Range("D2").Value = "4657"
$battingSheet.Range("D2").Style = $battingDefaultStyle

$battingSheet.Range("D3").NumberFormat = "@"
$battingSheet.Range("D3").Value = "4658"
$battingSheet.Range("D3").Style = $battingDefaultStyle

$battingSheet.Range("D4").NumberFormat = "@"
$battingSheet.Range("D4").Value = "4679"
$battingSheet.Range("D4").Style = $battingDefaultStyle

# --- 3. ODI Bowling sheet: MATCH_CARD_LINK -> MATCH_CODE --------------------

$bowlingSheet = $wb.Worksheets.Item("ODI Bowling")

$bowlingSheet.Range("B1").Value = "MATCH_CODE"

$bowlingDefaultStyle = $bowlingSheet.Range("A2").Style

$bowlingSheet.Range("B2").NumberFormat = "@"
$bowlingSheet.Range("B2").Value = "4657"
$bowlingSheet.Range("B2").Style = $bowlingDefaultStyle

$bowlingSheet.Range("B3").NumberFormat = "@"
$bowlingSheet.Range("B3").Value = "4658"
$bowlingSheet.Range("B3").Style = $bowlingDefaultStyle

$bowlingSheet.Range("B4").NumberFormat = "@"
$bowlingSheet.Range("B4").Value = "4679"
$bowlingSheet.Range("B4").Style = $bowlingDefaultStyle
